# Insert a new data row at row 60 (pushing existing rows 60-110 down to 61-111)
# and populate it with the new record's values, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 60:110 down by inserting a new row before row 60.
$ws.Rows.Item(60).Insert()

# New row 60 values (rest of the columns mirror the row that used to occupy
# row 60, i.e. "Sin especificar" variety and "$/saco 25 kilos" unit, which
# remain unchanged in the diff).
$ws.Cells.Item(60, 1).Value = 9
$ws.Cells.Item(60, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(60, 3).Value = "Metropolitana"
$ws.Cells.Item(60, 4).Value = 44651
$ws.Cells.Item(60, 4).NumberFormat = $ws.Cells.Item(61, 4).NumberFormat
$ws.Cells.Item(60, 5).Value = 13
$ws.Cells.Item(60, 6).Value = 100112022
$ws.Cells.Item(60, 7).Value = "Arveja Verde"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 15
$ws.Cells.Item(60, 11).Value = 32000
$ws.Cells.Item(60, 12).Value = 32000
$ws.Cells.Item(60, 13).Value = 32000
$ws.Cells.Item(60, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(60, 15).Value = "Carahue"
$ws.Cells.Item(60, 16).Value = 1280
$ws.Cells.Item(60, 17).Value = 25
$ws.Cells.Item(60, 18).Value = "Hortaliza"
